$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of master data (16th May refresh)
$newRows = @(
    @(10005, 110033, "eng", $true, "superadmin", "now()", "now()"),
    @(10005, 110034, "eng", $true, "superadmin", "now()", "now()"),
    @(10005, 110035, "eng", $true, "superadmin", "now()", "now()")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Move selection to the next empty row (matches post-entry Excel selection)
$ws.Range("A37:XFD1048576").Select()
